# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new value for column F }
$updates = @{
    "展览"   = @{ 2 = 318; 4 = 8288; 5 = 6051; 6 = 520; 7 = 103; 11 = 965 }
    "全部类型" = @{ 2 = 318; 4 = 8288; 5 = 6051; 6 = 520; 7 = 103; 15 = 965 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}

$wb.Save()
